$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.627.06'
$ws.Range('E2').Value = '  -0.04%  '

$ws.Range('D3').Value = '1.884.03'
$ws.Range('E3').Value = '  -0.15%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = $ws.Range('C4').Style
$ws.Range('E4').Value = '  -0.16%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.39'
$ws.Range('D5').Style = $ws.Range('C5').Style
$ws.Range('E5').Value = '  -0.27%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9995'
$ws.Range('D6').Style = $ws.Range('C6').Style
$ws.Range('E6').Value = '  -0.20%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4731'
$ws.Range('D7').Style = $ws.Range('C7').Style
$ws.Range('E7').Value = '  +0.08%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2896'
$ws.Range('D8').Style = $ws.Range('C8').Style
$ws.Range('E8').Value = '  -0.77%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06544'
$ws.Range('D9').Style = $ws.Range('C9').Style
$ws.Range('E9').Value = '  +0.46%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.41'
$ws.Range('D10').Style = $ws.Range('C10').Style
$ws.Range('E10').Value = '  -0.07%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '100.06'
$ws.Range('D11').Style = $ws.Range('C11').Style
$ws.Range('E11').Value = '  +3.76%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7626'
$ws.Range('D12').Style = $ws.Range('C12').Style
$ws.Range('E12').Value = '  +3.20%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07823'
$ws.Range('D13').Style = $ws.Range('C13').Style
$ws.Range('E13').Value = '  +0.47%  '

$ws.Range('D14').Value = '1.881.59'
$ws.Range('E14').Value = '  -0.57%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.242'
$ws.Range('D15').Style = $ws.Range('C15').Style
$ws.Range('E15').Value = '  +0.24%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '284.97'
$ws.Range('D16').Style = $ws.Range('C16').Style
$ws.Range('E16').Value = '  +0.73%  '

$ws.Range('D17').Value = '30.595.22'
$ws.Range('E17').Value = '  -0.15%  '

$ws.Range('E18').Value = '  -0.32%  '

$ws.Range('E19').Value = '  +0.35%  '

$ws.Range('E20').Value = '  -0.17%  '

$ws.Range('D21').Value = '2.125.90'
$ws.Range('E21').Value = '  -0.65%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.358'
$ws.Range('D22').Style = $ws.Range('C22').Style
$ws.Range('E22').Value = '  +0.98%  '

$ws.Range('E23').Value = '  -0.01%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.436'
$ws.Range('D24').Style = $ws.Range('C24').Style
$ws.Range('E24').Value = '  +3.01%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.182'
$ws.Range('D25').Style = $ws.Range('C25').Style
$ws.Range('E25').Value = '  -0.29%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.51'
$ws.Range('D26').Style = $ws.Range('C26').Style
$ws.Range('E26').Value = '  -0.56%  '

$ws.Range('E27').Value = '  +0.85%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.912'
$ws.Range('D28').Style = $ws.Range('C28').Style
$ws.Range('E28').Value = '  -0.01%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.09734'
$ws.Range('D29').Style = $ws.Range('C29').Style
$ws.Range('E29').Value = '  -0.14%  '

$ws.Range('E30').Value = '  -0.92%  '

$ws.Range('E31').Value = '  +0.97%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.256'
$ws.Range('D32').Style = $ws.Range('C32').Style

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.185'
$ws.Range('D33').Style = $ws.Range('C33').Style
$ws.Range('E33').Value = '  +0.24%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04850'
$ws.Range('D34').Style = $ws.Range('C34').Style
$ws.Range('E34').Value = '  -1.00%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.130'
$ws.Range('D35').Style = $ws.Range('C35').Style
$ws.Range('E35').Value = '  +0.43%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6993'
$ws.Range('D36').Style = $ws.Range('C36').Style
$ws.Range('E36').Value = '  +0.28%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.784'
$ws.Range('D37').Style = $ws.Range('C37').Style
$ws.Range('E37').Value = '  +2.32%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01909'
$ws.Range('D38').Style = $ws.Range('C38').Style
$ws.Range('E38').Value = '  +0.64%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.883'
$ws.Range('D39').Style = $ws.Range('C39').Style
$ws.Range('E39').Value = '  +1.67%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.316'
$ws.Range('D40').Style = $ws.Range('C40').Style
$ws.Range('E40').Value = '  +0.47%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '75.41'
$ws.Range('D41').Style = $ws.Range('C41').Style
$ws.Range('E41').Value = '  -0.94%  '

$ws.Range('E42').Value = '  -1.09%  '

$ws.Range('E43').Value = '  -0.50%  '

$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8404'
$ws.Range('D44').Style = $ws.Range('C44').Style
$ws.Range('E44').Value = '  +0.70%  '

$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9993'
$ws.Range('D45').Style = $ws.Range('C45').Style
$ws.Range('E45').Value = '  -0.22%  '

$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '101.46'
$ws.Range('D46').Style = $ws.Range('C46').Style
$ws.Range('E46').Value = '  -0.38%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.895'
$ws.Range('D47').Style = $ws.Range('C47').Style
$ws.Range('E47').Value = '  +3.93%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.037'
$ws.Range('D48').Style = $ws.Range('C48').Style
$ws.Range('E48').Value = '  +0.51%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '35.35'
$ws.Range('D49').Style = $ws.Range('C49').Style
$ws.Range('E49').Value = '  -0.37%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05787'
$ws.Range('D50').Style = $ws.Range('C50').Style
$ws.Range('E50').Value = '  +0.24%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3960'
$ws.Range('D51').Style = $ws.Range('C51').Style
$ws.Range('E51').Value = '  +0.11%  '
